$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 112
$ws.Range("H112").Value = 2031.8438
$ws.Range("I112").Value = 900.3333
$ws.Range("J112").Value = 2148.8965
$ws.Range("K112").Value = 2700.9999
$ws.Range("L112").Value = 6446.689499999999
$ws.Range("M112").Value = -1592.9999
$ws.Range("N112").Value = -8662.6895

# Row 138
$ws.Range("H138").Value = 2861535.2
$ws.Range("I138").Value = 5884592.5
$ws.Range("J138").Value = 6425.25
$ws.Range("K138").Value = 17653777.5
$ws.Range("L138").Value = 19275.75
$ws.Range("M138").Value = -17648637.5
$ws.Range("N138").Value = -29555.75

# Row 140
$ws.Range("H140").Value = 89080
$ws.Range("J140").Value = 88633
$ws.Range("L140").Value = 88633
$ws.Range("N140").Value = -98993


$ws = $wb.Worksheets.Item("ARM")
# Row 61
$ws.Range("H61").Value = 2106.6875
$ws.Range("I61").Value = 2075.7693
$ws.Range("J61").Value = 2240.6667
$ws.Range("K61").Value = 2075.7693
$ws.Range("L61").Value = 2240.6667
$ws.Range("M61").Value = -1863.7693
$ws.Range("N61").Value = -2664.6667

# Row 74
$ws.Range("H74").Value = 1343.6586
$ws.Range("I74").Value = 1009.2414
$ws.Range("K74").Value = 1009.2414
$ws.Range("M74").Value = -135.2414

# Row 77
$ws.Range("H77").Value = 1343.6586
$ws.Range("I77").Value = 1009.2414
$ws.Range("K77").Value = 5046.207
$ws.Range("M77").Value = -678.2070000000003

# Row 110
$ws.Range("H110").Value = 1890
$ws.Range("I110").Value = 1780
$ws.Range("J110").Value = 2000
$ws.Range("K110").Value = 1780
$ws.Range("L110").Value = 2000
$ws.Range("M110").Value = 265
$ws.Range("N110").Value = -6090

# Row 120
$ws.Range("H120").Value = 29903.166
$ws.Range("J120").Value = 29903.166
$ws.Range("L120").Value = 29903.166
$ws.Range("N120").Value = -39579.166

# Row 122
$ws.Range("H122").Value = 3247.3333
$ws.Range("I122").Value = 3215.75
$ws.Range("K122").Value = 9647.25
$ws.Range("M122").Value = -7197.25

# Row 123
$ws.Range("H123").Value = 58429
$ws.Range("J123").Value = 58429
$ws.Range("L123").Value = 58429
$ws.Range("N123").Value = -68229

# Row 134
$ws.Range("H134").Value = 67660
$ws.Range("J134").Value = 67660
$ws.Range("L134").Value = 67660
$ws.Range("N134").Value = -77800

# Row 136
$ws.Range("H136").Value = 2106.6875
$ws.Range("I136").Value = 2075.7693
$ws.Range("J136").Value = 2240.6667
$ws.Range("K136").Value = 6227.3079
$ws.Range("L136").Value = 6722.000100000001
$ws.Range("M136").Value = -3677.3079
$ws.Range("N136").Value = -11822.0001

# Row 140
$ws.Range("H140").Value = 48984.043
$ws.Range("J140").Value = 48984.043
$ws.Range("L140").Value = 48984.043
$ws.Range("N140").Value = -59344.043


$ws = $wb.Worksheets.Item("BSM")
# Row 22
$ws.Range("H22").Value = 12873.75
$ws.Range("I22").Value = 12873.75
$ws.Range("J22").Value = 0
$ws.Range("K22").Value = 12873.75
$ws.Range("L22").Value = 0
$ws.Range("M22").ClearContents()
$ws.Range("N22").Value = -12700.75

# Row 80
$ws.Range("H80").Value = 2315049.5
$ws.Range("I80").Value = 6173062
$ws.Range("J80").Value = 241.8
$ws.Range("K80").Value = 6173062
$ws.Range("L80").Value = 241.8
$ws.Range("M80").Value = -6172064
$ws.Range("N80").Value = -2237.8

# Row 83
$ws.Range("H83").Value = 2315049.5
$ws.Range("I83").Value = 6173062
$ws.Range("J83").Value = 241.8
$ws.Range("K83").Value = 30865310
$ws.Range("L83").Value = 1209
$ws.Range("M83").Value = -30860318
$ws.Range("N83").Value = -11193

# Row 132
$ws.Range("H132").Value = 73836.37
$ws.Range("J132").Value = 73836.37
$ws.Range("L132").Value = 73836.37
$ws.Range("N132").Value = -83956.37

# Row 134
$ws.Range("H134").Value = 457417.12
$ws.Range("I134").Value = 647556.3
$ws.Range("K134").Value = 1942668.9
$ws.Range("M134").Value = -1940133.9

# Row 140
$ws.Range("H140").Value = 49814.617
$ws.Range("J140").Value = 49814.617
$ws.Range("L140").Value = 49814.617
$ws.Range("N140").Value = -60174.617


$ws = $wb.Worksheets.Item("CRP")
# Row 31
$ws.Range("H31").Value = 3095.761
$ws.Range("I31").Value = 1856.28
$ws.Range("J31").Value = 4571.3335
$ws.Range("K31").Value = 1856.28
$ws.Range("L31").Value = 4571.3335
$ws.Range("M31").Value = -1561.28
$ws.Range("N31").Value = -5161.3335

# Row 34
$ws.Range("H34").Value = 3095.761
$ws.Range("I34").Value = 1856.28
$ws.Range("J34").Value = 4571.3335
$ws.Range("K34").Value = 1856.28
$ws.Range("L34").Value = 4571.3335
$ws.Range("M34").Value = -1654.28
$ws.Range("N34").Value = -4975.3335

# Row 134
$ws.Range("H134").Value = 1311.6774
$ws.Range("I134").Value = 1124.5555
$ws.Range("K134").Value = 3373.6665
$ws.Range("M134").Value = -838.6664999999998

# Row 135
$ws.Range("H135").Value = 94821.25
$ws.Range("J135").Value = 219523.33
$ws.Range("L135").Value = 219523.33
$ws.Range("N135").Value = -229663.33

# Row 138
$ws.Range("H138").Value = 59916.668
$ws.Range("J138").Value = 59916.668
$ws.Range("L138").Value = 59916.668
$ws.Range("N138").Value = -70196.66800000001


$ws = $wb.Worksheets.Item("CUL")
# Row 129
$ws.Range("H129").Value = 1924929.4
$ws.Range("I129").Value = 600
$ws.Range("J129").Value = 2175929
$ws.Range("K129").Value = 1800
$ws.Range("L129").Value = 6527787
$ws.Range("M129").Value = 3200
$ws.Range("N129").Value = -6537787


$ws = $wb.Worksheets.Item("GSM")
# Row 122
$ws.Range("H122").Value = 3781.9185
$ws.Range("I122").Value = 3651.697
$ws.Range("J122").Value = 4050.5
$ws.Range("K122").Value = 10955.091
$ws.Range("L122").Value = 12151.5
$ws.Range("M122").Value = -8505.091
$ws.Range("N122").Value = -17051.5

# Row 132
$ws.Range("H132").Value = 3442.6924
$ws.Range("I132").Value = 2218.8333
$ws.Range("J132").Value = 4491.7144
$ws.Range("K132").Value = 6656.499899999999
$ws.Range("L132").Value = 13475.1432
$ws.Range("M132").Value = -4126.499899999999
$ws.Range("N132").Value = -18535.1432

# Row 140
$ws.Range("H140").Value = 50250
$ws.Range("J140").Value = 50250
$ws.Range("L140").Value = 50250
$ws.Range("N140").Value = -60610


$ws = $wb.Worksheets.Item("LTW")
# Row 94
$ws.Range("H94").Value = 34000
$ws.Range("J94").Value = 34000
$ws.Range("L94").Value = 34000
$ws.Range("N94").Value = -35352

# Row 111
$ws.Range("H111").Value = 48715
$ws.Range("J111").Value = 48715
$ws.Range("L111").Value = 48715
$ws.Range("N111").Value = -56895

# Row 133
$ws.Range("H133").Value = 58306.668
$ws.Range("J133").Value = 58306.668
$ws.Range("L133").Value = 58306.668
$ws.Range("N133").Value = -63366.668

# Row 135
$ws.Range("H135").Value = 149727
$ws.Range("J135").Value = 149727
$ws.Range("L135").Value = 149727
$ws.Range("N135").Value = -159867

# Row 136
$ws.Range("H136").Value = 4462.476
$ws.Range("I136").Value = 4503.543
$ws.Range("K136").Value = 13510.629
$ws.Range("M136").Value = -10960.629

# Row 138
$ws.Range("H138").Value = 76960
$ws.Range("J138").Value = 76960
$ws.Range("L138").Value = 76960
$ws.Range("N138").Value = -87240


$ws = $wb.Worksheets.Item("WVR")
# Row 41
$ws.Range("H41").Value = 333337120
$ws.Range("I41").Value = 0
$ws.Range("J41").Value = 333337120
$ws.Range("K41").Value = 0
$ws.Range("L41").Value = 333337120
$ws.Range("M41").ClearContents()
$ws.Range("N41").Value = -333337900

# Row 132
$ws.Range("H132").Value = 1899.5588
$ws.Range("I132").Value = 1455.2222
$ws.Range("J132").Value = 3613.4285
$ws.Range("K132").Value = 4365.6666
$ws.Range("L132").Value = 10840.2855
$ws.Range("M132").Value = -1835.6666
$ws.Range("N132").Value = -15900.2855

# Row 139
$ws.Range("H139").Value = 60787.145
$ws.Range("J139").Value = 60787.145
$ws.Range("L139").Value = 60787.145
$ws.Range("N139").Value = -71067.14499999999

# Row 140
$ws.Range("H140").Value = 44500
$ws.Range("J140").Value = 44500
$ws.Range("L140").Value = 44500
$ws.Range("N140").Value = -54860

